$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) holds text-formatted numbers (e.g. thousands separators
# using '.' and values like 0.00000000108). Force text format first so
# Excel doesn't silently reinterpret the new values as real numbers/dates.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '27.945.74'
$ws.Range('D3').Value = '1.739.53'
$ws.Range('E3').Value = '  -4.69%  '
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  -0.34%  '
$ws.Range('D5').Value = '226.37'
$ws.Range('E5').Value = '  -3.74%  '
$ws.Range('D6').Value = '0.5788'
$ws.Range('E6').Value = '  -3.37%  '
$ws.Range('E7').Value = '  -0.24%  '
$ws.Range('D8').Value = '0.2736'
$ws.Range('E8').Value = '  -0.91%  '
$ws.Range('D9').Value = '23.18'
$ws.Range('E9').Value = '  -1.10%  '
$ws.Range('D10').Value = '0.06617'
$ws.Range('E10').Value = '  -4.69%  '
$ws.Range('D11').Value = '0.07542'
$ws.Range('E11').Value = '  -0.88%  '
$ws.Range('D12').Value = '1.741.65'
$ws.Range('E12').Value = '  -5.83%  '
$ws.Range('D13').Value = '4.708'
$ws.Range('E13').Value = '  -0.49%  '
$ws.Range('D14').Value = '0.6023'
$ws.Range('E14').Value = '  -4.16%  '
$ws.Range('D15').Value = '1.976.14'
$ws.Range('E15').Value = '  -4.70%  '
$ws.Range('D16').Value = '74.70'
$ws.Range('E16').Value = '  -3.53%  '
$ws.Range('D17').Value = '0.000008760'
$ws.Range('E17').Value = '  -10.77%  '
$ws.Range('D18').Value = '27.930.67'
$ws.Range('E18').Value = '  -3.71%  '
$ws.Range('D19').Value = '5.317'
$ws.Range('E19').Value = '  -4.00%  '
$ws.Range('D20').Value = '1.001'
$ws.Range('E20').Value = '  -0.26%  '
$ws.Range('D21').Value = '205.53'
$ws.Range('E21').Value = '  -4.74%  '
$ws.Range('D22').Value = '11.29'
$ws.Range('E22').Value = '  -2.19%  '
$ws.Range('D23').Value = '6.626'
$ws.Range('E23').Value = '  -3.30%  '
$ws.Range('E24').Value = '  -0.26%  '
$ws.Range('D25').Value = '150.22'
$ws.Range('E25').Value = '  -3.53%  '
$ws.Range('D26').Value = '8.050'
$ws.Range('E26').Value = '  +1.35%  '
$ws.Range('D27').Value = '0.1233'
$ws.Range('E27').Value = '  -4.14%  '
$ws.Range('D28').Value = '16.18'
$ws.Range('E28').Value = '  -1.91%  '
$ws.Range('E29').Value = '  -2.09%  '
$ws.Range('D30').Value = '0.06138'
$ws.Range('E30').Value = '  -4.32%  '
$ws.Range('D31').Value = '1.393'
$ws.Range('E31').Value = '  -3.35%  '
$ws.Range('D32').Value = '3.741'
$ws.Range('E32').Value = '  -2.03%  '
$ws.Range('D33').Value = '3.734'
$ws.Range('E33').Value = '  -1.18%  '
$ws.Range('D34').Value = '1.672'
$ws.Range('E34').Value = '  -2.90%  '
$ws.Range('D35').Value = '1.036'
$ws.Range('E35').Value = '  -5.26%  '
$ws.Range('D36').Value = '0.6412'
$ws.Range('E36').Value = '  -0.67%  '
$ws.Range('D37').Value = '2.417'
$ws.Range('E37').Value = '  -4.93%  '
$ws.Range('D38').Value = '2.716'
$ws.Range('E38').Value = '  -1.44%  '
$ws.Range('D39').Value = '0.01667'
$ws.Range('E39').Value = '  -4.88%  '
$ws.Range('D40').Value = '1.124.64'
$ws.Range('E40').Value = '  -0.78%  '
$ws.Range('D41').Value = '6.151'
$ws.Range('E41').Value = '  -6.66%  '
$ws.Range('D42').Value = '0.8759'
$ws.Range('E42').Value = '  -1.83%  '
$ws.Range('D43').Value = '1.004'
$ws.Range('E43').Value = '  -0.01%  '
$ws.Range('D44').Value = '99.90'
$ws.Range('E44').Value = '  -0.99%  '
$ws.Range('D45').Value = '1.888.77'
$ws.Range('E45').Value = '  -4.97%  '
$ws.Range('D46').Value = '59.32'
$ws.Range('E46').Value = '  -4.47%  '
$ws.Range('B47').Value = 'BabyDogeCoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D47').Value = '0.00000000108'
$ws.Range('E47').Value = '  -4.27%  '
$ws.Range('B48').Value = 'RenderToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D48').Value = '1.577'
$ws.Range('E48').Value = '  -2.28%  '
$ws.Range('D49').Value = '8.261'
$ws.Range('E49').Value = '  -1.99%  '
$ws.Range('E50').Value = '  -2.35%  '
$ws.Range('B51').Value = 'Mantle'
$ws.Range('C51').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D51').Value = '0.4412'
$ws.Range('E51').Value = '  -2.76%  '
